$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.649.32'
$ws.Range('D3').Value = '1.611.02'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = "'212.45"
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').Value = "'0.520"
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('D7').Value = "'0.993"
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').Value = "'28.90"
$ws.Range('E8').Value = '  +7.64%  '
$ws.Range('D9').Value = "'0.257"
$ws.Range('E9').Value = '  +2.92%  '
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').Value = "'0.0906"
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = '1.842.07'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '1.610.25'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = "'0.569"
$ws.Range('E14').Value = '  +6.11%  '
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Value = '29.656.00'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = "'8.68"
$ws.Range('E17').Value = '  +14.05%  '
$ws.Range('D18').Value = "'64.62"
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').Value = "'240.60"
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').Value = '0.0₃0704'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('E23').Value = '  +4.83%  '
$ws.Range('D24').Value = "'2.11"
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('D25').Value = "'156.51"
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('D26').Value = "'15.61"
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  +2.48%  '
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('D34').Value = '1.437.77'
$ws.Range('E34').Value = '  +1.67%  '
$ws.Range('E35').Value = '  +5.87%  '
$ws.Range('E36').Value = '  +1.78%  '
$ws.Range('E37').Value = '  +2.93%  '
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('E40').Value = '  +3.78%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = "'0.0506"
$ws.Range('E41').Value = '  +6.74%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = "'2.00"
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').Value = "'0.823"
$ws.Range('E43').Value = '  +3.90%  '
$ws.Range('D44').Value = "'54.04"
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').Value = "'69.60"
$ws.Range('E45').Value = '  +6.02%  '
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('E47').Value = '  +20.48%  '
$ws.Range('D48').Value = "'5.44"
$ws.Range('E48').Value = '  +3.18%  '
$ws.Range('D49').Value = '1.750.80'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('E51').Value = '  -1.04%  '
